$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.295175552368164
$ws.Range("B1").Value = 2.507584571838379
$ws.Range("C1").Value = 1.273186445236206
$ws.Range("D1").Value = 0.5792209506034851
$ws.Range("E1").Value = 0.4919843375682831
